$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Стятюгин Александр): fill in the previously-missing homework
# scores (ДЗ_1, ДЗ_2, ДЗ_3, ДЗ_4, лаб_1, ДЗ_7) with 5 — Сумма (L16)
# recalculates automatically via its shared SUM formula.
$ws.Range("C16:F16").Value = 5
$ws.Range("I16:J16").Value = 5

# Add a note in O16 documenting the change.
$ws.Range("O16").Value = "выаы"

# Move the active selection to the cell that was just edited.
$ws.Range("O16").Select()
